# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New account-statement data (rows 16-27): Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico
$rows = @(
    @("CC", "1049827115", "JORGE LUIS AYOLA DIAZ", "2305", 32000, 1160000),
    @("CC", "1049827115", "JORGE LUIS AYOLA DIAZ", "2212", 40000, 1160000),
    @("CC", "1049827115", "JORGE LUIS AYOLA DIAZ", "2211", 40000, 1160000),
    @("CC", "1047415271", "YIRA TATIANA DIAZ NOVOA", "2305", 32000, 1319997),
    @("CC", "1047415271", "YIRA TATIANA DIAZ NOVOA", "2212", 40000, 1319997),
    @("CC", "1047415271", "YIRA TATIANA DIAZ NOVOA", "2211", 40000, 1319997),
    @("CC", "73166479", "MAXIMILIANO MARRUGO POLO", "2305", 37120, 1160000),
    @("CC", "73560811", "ROMAN CARRIAZO GOENAGA", "2305", 37120, 1160000),
    @("CC", "1051886209", "JORGE ENRIQUE DEULOFEUTT CASTILLO", "2305", 37120, 1160000),
    @("CC", "73106277", "FIDEL CASTILLO GODOY", "2305", 37120, 1160000),
    @("CC", "1049824118", "ALFONSO MONTALBAN MERCADO", "2305", 37120, 1160000),
    @("CC", "1051891144", "LUIS FERNANDO BATISTA MORALES", "2305", 37120, 1160000)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]
    $ws.Cells.Item($r, 6).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[5]
}
